# ROP_FSB.pptx - slide 16 ("What is Format String?")
# The single sentence in the content placeholder is reworded: the two
# middle words "지정 서식자가" / "서식 지정자가" are swapped, and the
# sentence is additionally split into three separate runs (matching the
# way PowerPoint leaves behind run boundaries after an in-place edit of
# a selected sub-string).
#
# Before (1 run):
#   "아래와 같이 다양한 지정 서식자가 존재함" + "." (separate run, unchanged)
# After (3 runs):
#   "아래와 같이 " | "다양한 서식 지정자가 " | "존재함" + "." (unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldSentence = "아래와 같이 다양한 지정 서식자가 존재함"
$newSentence = "아래와 같이 다양한 서식 지정자가 존재함"

if ($tr.Text.IndexOf($oldSentence) -ne 0) {
    throw "Expected slide 16 / shape 2 to start with the known sentence; found: $($tr.Text)"
}

$lead = "아래와 같이 "
$mid  = "다양한 서식 지정자가 "
$tail = "존재함"

# Locate the run holding the Korean sentence (it is the first characters
# of the paragraph) and rewrite its text in place (word-order swap).
$target = $tr.Characters(1, $oldSentence.Length)
$target.Text = $newSentence

# Re-split the (now updated) sentence into three runs along the same
# boundaries shown in the diff, so each chunk becomes its own <a:r>.
$part1 = $tr.Characters(1, $lead.Length)
$part1.Text = $part1.Text

$part2 = $tr.Characters($lead.Length + 1, $mid.Length)
$part2.Text = $part2.Text

Write-Host "Final text:" $tr.Text
